$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "五洲新春"
$ws.Range("B2").Value = "协鑫集成"
$ws.Range("C2").Value = "协鑫集成"
$ws.Range("A3").Value = "数据港"
$ws.Range("B3").Value = "数据港"
$ws.Range("C3").Value = "巨力索具"
$ws.Range("A4").Value = "协鑫集成"
$ws.Range("B4").Value = "杭电股份"
$ws.Range("C4").Value = "神剑股份"
$ws.Range("A5").Value = "利欧股份"
$ws.Range("B5").Value = "神剑股份"
$ws.Range("C5").Value = "杭电股份"
$ws.Range("A6").Value = "杭电股份"
$ws.Range("B6").Value = "银河电子"
$ws.Range("C6").Value = "航天发展"
$ws.Range("A7").Value = "神剑股份"
$ws.Range("B7").Value = "洲际油气"
$ws.Range("C7").Value = "利欧股份"
$ws.Range("A8").Value = "巨力索具"
$ws.Range("B8").Value = "巨力索具"
$ws.Range("C8").Value = "浙文互联"
$ws.Range("A9").Value = "银河电子"
$ws.Range("B9").Value = "五洲新春"
$ws.Range("C9").Value = "洲际油气"
$ws.Range("A10").Value = "洲际油气"
$ws.Range("B10").Value = "协鑫能科"
$ws.Range("C10").Value = "银河电子"
$ws.Range("A11").Value = "科瑞技术"
$ws.Range("B11").Value = "湖南黄金"
$ws.Range("C11").Value = "五洲新春"
$ws.Range("A12").Value = "天奇股份"
$ws.Range("B12").Value = "杉杉股份"
$ws.Range("C12").Value = "天奇股份"
$ws.Range("A13").Value = "长飞光纤"
$ws.Range("B13").Value = "利欧股份"
$ws.Range("C13").Value = "数据港"
$ws.Range("A14").Value = "航天发展"
$ws.Range("B14").Value = "科瑞技术"
$ws.Range("C14").Value = "雷科防务"
$ws.Range("A15").Value = "雷科防务"
$ws.Range("B15").Value = "天奇股份"
$ws.Range("C15").Value = "中超控股"
$ws.Range("A16").Value = "湖南黄金"
$ws.Range("B16").Value = "贵州茅台"
$ws.Range("C16").Value = "白银有色"
$ws.Range("A17").Value = "浙文互联"
$ws.Range("B17").Value = "三变科技"
$ws.Range("C17").Value = "三变科技"
$ws.Range("A18").Value = "三变科技"
$ws.Range("B18").Value = "雷科防务"
$ws.Range("C18").Value = "湖南白银"
$ws.Range("A19").Value = "福莱新材"
$ws.Range("B19").Value = "东方财富"
$ws.Range("C19").Value = "TCL中环"
$ws.Range("A20").Value = "智立方"
$ws.Range("B20").Value = "百川股份"
$ws.Range("C20").Value = "平潭发展"
$ws.Range("A21").Value = "协鑫能科"
$ws.Range("B21").Value = "浙文互联"
$ws.Range("C21").Value = "锋龙股份"
